$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.902.50"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.875.20"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("E4").Value = "  -0.28%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.01"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  -1.63%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  -0.22%  "
$__style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5310"
$ws.Range("D7").Style = $__style
$ws.Range("E7").Value = "  +1.52%  "
$__style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3756"
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = "  -1.19%  "
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07158"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  -1.72%  "
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.65"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = "  +1.60%  "
$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8859"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = "  -2.14%  "
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08141"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").Value = "1.922.23"
$ws.Range("E13").Value = "  +1.40%  "
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.10"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = "  -2.49%  "
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.279"
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = "  -1.32%  "
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = "  -0.27%  "
$__style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.72"
$ws.Range("D17").Style = $__style
$ws.Range("E17").Value = "  +0.39%  "
$__style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008553"
$ws.Range("D18").Style = $__style
$ws.Range("E18").Value = "  -1.15%  "
$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "26.997.89"
$ws.Range("E20").Value = "  -0.64%  "
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.980"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("E22").Value = "  -0.80%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.383"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  -1.29%  "
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.46"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = "  -1.44%  "
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.275"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("E26").Value = "  -0.37%  "
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.99"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = "  -1.51%  "
$__style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.53"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("E29").Value = "  -1.75%  "
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.577"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = "  -5.90%  "
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09102"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  -1.44%  "
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7996"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = "  +0.85%  "
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04987"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = "  -1.18%  "
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.173"
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = "  -4.08%  "
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.988"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("B36").Value = "TheSandbox"
$ws.Range("C36").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5868"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  +2.34%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.206"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  -5.28%  "
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.605"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  -2.03%  "
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01952"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = "  -2.18%  "
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.610"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.900"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = "  -1.34%  "
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.28"
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = "  +0.01%  "
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5064"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +3.40%  "
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1495"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  -1.44%  "
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9995"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = "  -0.34%  "
$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.01"
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = "  -1.44%  "
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.611"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = "  -1.83%  "
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.98"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = "  -1.58%  "
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06030"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  +1.24%  "
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.54"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = "  -2.37%  "
